$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 77

# Column A holds the date as plain text (matching the existing rows, which
# are shared strings rather than real Excel dates). Writing the literal
# string directly would make Excel auto-convert it into a date serial
# number, so instead we compute it as a text formula and then flatten the
# formula down to its static (text) value via Copy + PasteSpecial values.
# That preserves the "text" cell type without leaving any stray cell
# formatting/style behind.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.Formula = '="2020-08-15"'
$cellA.Copy()
$cellA.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 2).Value = 517714
$ws.Cells.Item($newRow, 3).Value = 568359
$ws.Cells.Item($newRow, 4).Value = 84934
$ws.Cells.Item($newRow, 5).Value = 56543
$ws.Cells.Item($newRow, 6).Value = 26.25
